# Kết quả test Với đầu vào có và ko có Volumes
# Add two new result tables to the "ANN" sheet:
#   - G1:K8   "ANN.NET Volume"   (same layout as the existing B1:E8 table)
#   - A10:E17 "ANN.NET NoVolume" (same layout as the existing B1:E8 table,
#                                 shifted down two rows)
# Also clears out the now-unused column F on "K-SVMeans" and touches the
# selection/active-sheet state on every sheet so it matches the author's
# final click-state (ANN ends up the active/selected tab).

$wb  = $excel.ActiveWorkbook
$ann = $wb.Worksheets.Item("ANN")
$svm = $wb.Worksheets.Item("SVM")
$ksv = $wb.Worksheets.Item("K-SVMeans")

# ---------------------------------------------------------------------
# ANN sheet: new "ANN.NET Volume" table at G1:K8 (mirrors B1:E8 / H column
# offset by 6 from B).
# ---------------------------------------------------------------------

$ann.Range("G1:K1").Merge()
$ann.Range("G1").Value = "ANN.NET Volume"
$ann.Range("G1:K1").HorizontalAlignment = -4108  # xlCenter

$ann.Range("H2").Value = "period = 1"
$ann.Range("I2").Value = "period = 5"
$ann.Range("J2").Value = "period = 10"
$ann.Range("K2").Value = "period = 30"

$ann.Range("G3").Value = "BT6"
$ann.Range("H3").Value = 68.1
$ann.Range("I3").Value = 67.78
$ann.Range("J3").Value = 59.33
$ann.Range("K3").Value = 41.3

$ann.Range("G4").Value = "DHG"
$ann.Range("H4").Value = 56.98
$ann.Range("I4").Value = 56.75
$ann.Range("J4").Value = 56.52
$ann.Range("K4").Value = 47.22

$ann.Range("G5").Value = "FPT"
$ann.Range("H5").Value = 67.02
$ann.Range("I5").Value = 32.62
$ann.Range("J5").Value = 43.01
$ann.Range("K5").Value = 47.8

$ann.Range("G6").Value = "VIS"
$ann.Range("H6").Value = 42.16
$ann.Range("I6").Value = 41.84
$ann.Range("J6").Value = 51.91
$ann.Range("K6").Value = 55.86

$ann.Range("G7").Value = "VNM"
$ann.Range("H7").Value = 59.05
$ann.Range("I7").Value = 56.7
$ann.Range("J7").Value = 54.34
$ann.Range("K7").Value = 41.15

$ann.Range("G8").Value = "Total"
$ann.Range("H8").Formula = "=AVERAGE(H3:H7)"
$ann.Range("I8:K8").Formula = "=AVERAGE(I3:I7)"
$ann.Range("H8").Font.Color = 255
$ann.Range("I8").Font.Color = 255

# widen column K a little (closest reachable width to the author's 11.29)
$ann.Columns.Item(11).ColumnWidth = 10.45

# ---------------------------------------------------------------------
# ANN sheet: new "ANN.NET NoVolume" table at A10:E17 (same shape as the
# original B1:E8 table, two rows further down).
# ---------------------------------------------------------------------

$ann.Range("A10:E10").Merge()
$ann.Range("A10").Value = "ANN.NET NoVolume"
$ann.Range("A10:E10").HorizontalAlignment = -4108  # xlCenter

$ann.Range("B11").Value = "period = 1"
$ann.Range("C11").Value = "period = 5"
$ann.Range("D11").Value = "period = 10"
$ann.Range("E11").Value = "period = 30"

$ann.Range("A12").Value = "BT6"
$ann.Range("B12").Value = 72.85
$ann.Range("C12").Value = 66.1
$ann.Range("D12").Value = 44.01
$ann.Range("E12").Value = 58.45

$ann.Range("A13").Value = "DHG"
$ann.Range("B13").Value = 56.98
$ann.Range("C13").Value = 56.75
$ann.Range("D13").Value = 44.02
$ann.Range("E13").Value = 31.66

$ann.Range("A14").Value = "FPT"
$ann.Range("B14").Value = 43.61
$ann.Range("C14").Value = 60.96
$ann.Range("D14").Value = 67.2
$ann.Range("E14").Value = 32.41

$ann.Range("A15").Value = "VIS"
$ann.Range("B15").Value = 49.72
$ann.Range("C15").Value = 38.58
$ann.Range("D15").Value = 60.65
$ann.Range("E15").Value = 52.51

$ann.Range("A16").Value = "VNM"
$ann.Range("B16").Value = 56.46
$ann.Range("C16").Value = 57.14
$ann.Range("D16").Value = 50
$ann.Range("E16").Value = 37.61

$ann.Range("A17").Value = "Total"
$ann.Range("B17").Formula = "=AVERAGE(B12:B16)"
$ann.Range("C17:E17").Formula = "=AVERAGE(C12:C16)"
$ann.Range("B17").Font.Color = 255
$ann.Range("C17").Font.Color = 255

# ---------------------------------------------------------------------
# K-SVMeans sheet: the "No Volume" column F is gone from this result set.
# ---------------------------------------------------------------------

$ksv.Range("F4:F9").Clear()

# ---------------------------------------------------------------------
# Selection / active-sheet bookkeeping, in click order: SVM, K-SVMeans,
# then ANN last so ANN ends up the active tab (matches the commit).
# ---------------------------------------------------------------------

$svm.Range("B9").Select()
$ksv.Range("B8").Select()
$ann.Range("F16").Select()
